# Summary.xlsx — add a new "Test 4" run and roll the reference in "Test 3"
# forward to compare against "Test 2" instead of "Test 1".

$wb = $excel.ActiveWorkbook

# --- Add the new "Test 4" sheet after "Test 3" by cloning it (keeps layout/format) ---
$ws3 = $wb.Worksheets.Item("Test 3")
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "Test 4"

# --- "Test 3": only the Delta formula's denominator changes (Test 1 -> Test 2) ---
$ws3.Range("B6").Formula = "=('Test 2'!B5-'Test 3'!B5)/'Test 2'!B5"

# --- "Test 4": drop in this run's raw timestamps (B2:B5 formulas already cloned) ---
$ws4.Range("A1").Value = 1537146896653
$ws4.Range("A2").Value = 1537154785973
$ws4.Range("A3").Value = 1537162973182
$ws4.Range("A4").Value = 1537170674441

# The Delta row now compares Test 3 against the new Test 4
$ws4.Range("B6").Formula = "=('Test 3'!B5-'Test 4'!B5)/'Test 3'!B5"

[void]$ws4.Range("B7").Select()

# "Test 4" is now the active tab
$ws4.Activate()
